$wb = $excel.ActiveWorkbook

$wsZone = $wb.Worksheets.Item("ZONE_CALCULATORS")
$wsZone.Range("B5").Value = 156
$wsZone.Range("B6").Value = 3496
$wsZone.Range("B7").Value = 219
$wsZone.Range("B41").Value = 1016
$wsZone.Range("B42").Value = 71
$wsZone.Range("B76").Value = 0
$wsZone.Range("B77").Value = 0
$wsZone.Range("B112").Value = 0
$wsZone.Range("B147").Value = 0

$wsRes = $wb.Worksheets.Item("RESOURCE_MGR")
$wsRes.Range("C27").Value = 156
$wsRes.Range("B38").Value = 156
